$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '45.979.82'
$ws.Range('E2').Value = '  -2.38%  '

Set-TextValue $ws.Range('D3') '2.662.12'
$ws.Range('E3').Value = '  +0.38%  '

$ws.Range('E4').Value = '  -0.08%  '

Set-TextValue $ws.Range('D5') '310.38'
$ws.Range('E5').Value = '  -1.68%  '

Set-TextValue $ws.Range('D6') '98.60'
$ws.Range('E6').Value = '  -5.75%  '

$ws.Range('E7').Value = '  -2.32%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('E9').Value = '  -3.15%  '

Set-TextValue $ws.Range('D10') '38.61'
$ws.Range('E10').Value = '  -2.66%  '

Set-TextValue $ws.Range('D11') '0.0848'
$ws.Range('E11').Value = '  -1.51%  '

$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D12') '8.09'
$ws.Range('E12').Value = '  -4.45%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D13') '3.065.17'
$ws.Range('E13').Value = '  +0.25%  '

$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D14') '0.108'
$ws.Range('E14').Value = '  +0.65%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D15') '2.653.72'
$ws.Range('E15').Value = '  -0.79%  '

$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D16') '0.935'
$ws.Range('E16').Value = '  -1.29%  '

$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D17') '15.13'
$ws.Range('E17').Value = '  -1.76%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range('D18') '45.994.17'
$ws.Range('E18').Value = '  -3.54%  '

$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D19') '0.0000102'
$ws.Range('E19').Value = '  -2.33%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D20') '6.84'
$ws.Range('E20').Value = '  -0.48%  '

$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D21') '12.88'
$ws.Range('E21').Value = '  -4.66%  '

$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D22') '74.87'
$ws.Range('E22').Value = '  +2.61%  '

$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range('D23') '284.47'
$ws.Range('E23').Value = '  +2.80%  '

$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D24') '3.07'
$ws.Range('E24').Value = '  -1.95%  '

$ws.Range('B25').Value = 'EthereumClassic'
$ws.Range('C25').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D25') '31.53'
$ws.Range('E25').Value = '  +0.32%  '

Set-TextValue $ws.Range('D26') '2.23'
$ws.Range('E26').Value = '  -0.24%  '

$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D27') '1.00'
$ws.Range('E27').Value = '  +0.09%  '

$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D28') '10.59'
$ws.Range('E28').Value = '  -2.31%  '

$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D29') '2.24'
$ws.Range('E29').Value = '  -3.05%  '

Set-TextValue $ws.Range('D30') '38.68'
$ws.Range('E30').Value = '  -7.43%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D31') '6.28'
$ws.Range('E31').Value = '  -0.03%  '

$ws.Range('B32').Value = 'LidoDAOToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D32') '3.81'
$ws.Range('E32').Value = '  -1.05%  '

$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D33') '2.39'
$ws.Range('E33').Value = '  +2.47%  '

$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D34') '156.26'
$ws.Range('E34').Value = '  +1.63%  '

$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D35') '0.0842'
$ws.Range('E35').Value = '  -1.53%  '

$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D36') '2.82'
$ws.Range('E36').Value = '  -1.96%  '

$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D37') '0.123'
$ws.Range('E37').Value = '  +1.08%  '

$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D38') '25.84'
$ws.Range('E38').Value = '  +5.63%  '

$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D39') '0.124'
$ws.Range('E39').Value = '  -0.61%  '

$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range('D40') '15.83'
$ws.Range('E40').Value = '  -6.54%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D41') '0.0329'
$ws.Range('E41').Value = '  -2.05%  '

$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D42') '3.59'
$ws.Range('E42').Value = '  -4.08%  '

$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D43') '2.175.56'
$ws.Range('E43').Value = '  +4.25%  '

Set-TextValue $ws.Range('D44') '3.97'
$ws.Range('E44').Value = '  -8.74%  '

$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D45') '0.998'
$ws.Range('E45').Value = '  -0.07%  '

$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue $ws.Range('D46') '94.44'
$ws.Range('E46').Value = '  -4.50%  '

$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D47') '112.18'
$ws.Range('E47').Value = '  -3.73%  '

Set-TextValue $ws.Range('D48') '9.29'
$ws.Range('E48').Value = '  -0.11%  '

$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range('D49') '2.914.45'
$ws.Range('E49').Value = '  +0.11%  '

$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D50') '0.201'
$ws.Range('E50').Value = '  -2.33%  '

$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D51') '1.74'
$ws.Range('E51').Value = '  -7.88%  '
